$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row above the current row 2 (LOZANO MOLINA TITO JERSON),
# shifting the existing advisors (and the totals row) down by one. This
# makes room for the new advisor "AVILA TORRES RAFAEL ALEJANDRO".
$ws1.Rows.Item(2).Insert()

# Excel's row-insert carries the formatting of the row above (the bold,
# centered header row) into the freshly inserted row. Reset it back to
# plain/default formatting before applying the real values, matching the
# other data rows.
$ws1.Rows.Item(2).ClearFormats()

$ws1.Cells.Item(2, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(2, 2).Value = "AVILA TORRES RAFAEL ALEJANDRO"
for ($col = 3; $col -le 14; $col++) {
    $ws1.Cells.Item(2, $col).Value = 0
    $ws1.Cells.Item(2, $col).NumberFormat = '"$"#,##0.00'
}

# Update the "0 de 2" totals row (now row 5) to reflect 3 advisors.
for ($col = 3; $col -le 14; $col++) {
    $ws1.Cells.Item(5, $col).Value = "0 de 3"
}

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(2).Insert()
$ws2.Rows.Item(2).ClearFormats()

$ws2.Cells.Item(2, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(2, 2).Value = "AVILA TORRES RAFAEL ALEJANDRO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(2, $col).Value = 0
    $ws2.Cells.Item(2, $col).NumberFormat = '"$"#,##0.00'
}
